$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3:H18").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G14").Value = 1
